$wb = $excel.ActiveWorkbook

# --- Sheet "Datasets" (sheet1): two new rows of dataset metadata ---
$ws1 = $wb.Worksheets.Item("Datasets")

# Fill row 10 before row 9 so the shared-strings table is built up in the
# same order the source workbook uses.
$ws1.Range("A10").Value = "merged-files-final-selected-features-2023-12-12"
$ws1.Range("D10").Value = "Automatic IQR OD"
$ws1.Range("B10").Value = 2804
$ws1.Range("C10").Value = 6

$ws1.Range("A9").Value = "merged-files-final-2023-12-12"
$ws1.Range("D9").Value = "Automatic IQR OD"
$ws1.Range("B9").Value = 2804
$ws1.Range("C9").Value = 17

# New rows are highlighted in bold green, matching the rest of the sheet's
# direct-formatting conventions.
$newDataRows = $ws1.Range("A9:D10")
$newDataRows.Font.Bold = $true
$newDataRows.Font.Color = 5287936

# --- Sheet "Measurements" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Measurements")

# C3 / C6 drop their stale per-cell formatting override and fall back to
# column C's default (left/top aligned) look.
foreach ($addr in @("C3", "C6")) {
    $cell = $ws2.Range($addr)
    $cell.HorizontalAlignment = -4131
    $cell.VerticalAlignment = -4160
}

# Two new measurement log rows.
$ws2.Range("A8").Value = "12/12/2023"
$ws2.Range("B8").Value = "all"
$ws2.Range("C8").Value = "merged-files-final-selected-features-2023-12-12"
$ws2.Range("D8").Value = "Calculated on Linux Home"
$ws2.Rows.Item(8).RowHeight = 17

$ws2.Range("A9").Value = "12/13/2023"
$ws2.Range("B9").Value = "all"
$ws2.Range("C9").Value = "merged-files-final-selected-features-2023-12-12"
$ws2.Range("D9").Value = "Calculated on Linux ZZ"
$ws2.Rows.Item(9).RowHeight = 17

# --- Selection bookkeeping: keep "Measurements" the active tab, but leave
# each sheet's own cursor on the last edited cell. ---
$ws1.Activate() | Out-Null
$ws1.Range("A10").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("D9").Select() | Out-Null
